$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Writes a literal text value into a cell without letting Excel's
    # autodetection turn date-looking strings (e.g. "2020.04.23") into
    # real dates. We do this by entering it as a formula that evaluates
    # to the literal string, then converting the cell to a static value
    # via copy / paste-special-values (this avoids creating any new,
    # unused cell style in the process).
    $r = $ws.Range($range)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163) # xlPasteValues
    $excel.CutCopyMode = 0
}

# Insert two blank columns before column B; this shifts the existing
# "hours" column (B) to D and the "what" column (C) to E, matching the
# new layout used to add start/end-time tracking in columns B and C.
$ws.Columns("B:C").Insert()

# ---- Row 2 (2020.02.08): hours 4 -> time-of-day duration 4:00 ----
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("D2").Value = [double](4/24)

# ---- Row 3 (2020.02.22): hours 1.5 -> time-of-day duration 1:30 ----
$ws.Range("D3").NumberFormat = "h:mm"
$ws.Range("D3").Value = [double](1.5/24)

# ---- Row 4: new entry, 1 hour, "meeting with Konrad" ----
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("D4").Value = [double](1/24)
Set-TextValue "E4" "meeting with Konrad"

# ---- Row 5: new entry, dated 2020.04.23, 1 hour, "meeting with Konrad" ----
Set-TextValue "A5" "2020.04.23"
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("D5").Value = [double](1/24)
Set-TextValue "E5" "meeting with Konrad"

# ---- Row 6: new entry, dated 2020.05.04, tracked via start/end time ----
Set-TextValue "A6" "2020.05.04"
$ws.Range("B6").NumberFormat = "h:mm"
$ws.Range("B6").Value = 0.5
$ws.Range("C6").NumberFormat = "h:mm"
$ws.Range("C6").Value = 0.54166666666666663
$ws.Range("D6").NumberFormat = "h:mm"
$ws.Range("D6").Formula = "=C6-B6"

# ---- Row 32: total hours logged ----
Set-TextValue "A32" "total"
$ws.Range("D32").NumberFormat = "h:mm"
$ws.Range("D32").Formula = "=SUM(D2:D31)"

# Restore the selection to match the author's last-worked-on cell.
$ws.Range("E16").Select()
